# Insert a new data row at row 2 (pushing the existing rows down by one),
# then fill in the new company's details. Clear any formatting the insert
# operation may have copied in so the new row matches the plain (unstyled)
# look of the other data rows, then restore column A's header-like style
# (bold, bordered, centered) by copying it from the row below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown, [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromRightOrBelow)
$ws.Range("A2:Q2").ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "南京伯索网络科技有限公司（PLASO）"
$ws.Range("C2").Value = "秦淮区"
$ws.Range("F2").Value = "9:00-18:00"
$ws.Range("G2").Value = "1h"
$ws.Range("H2").Value = "124 加班，35 正常；大小周"
$ws.Range("I2").Value = "基数南京底薪，比例 8%"
$ws.Range("J2").Value = "一般无"
$ws.Range("K2").Value = "3个月8折"
$ws.Range("L2").Value = "网吧工位"
$ws.Range("M2").Value = "入职一年后才有，每年加一天"
$ws.Range("N2").Value = "企业微信打卡，每月三次迟到机会"
$ws.Range("Q2").Value = "2022-02-06 13:26:16"

# Renumber column A (serial index) for all the rows that shifted down.
For ($r = 3; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
